$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 23 (SM / Biomass Taken Home (kg)) stats: remove outlier, recompute summary stats
$ws.Range("C23").Value = 17
$ws.Range("D23").Value = 0.5145744069746477
$ws.Range("E23").Value = 0.3493109811986957
$ws.Range("F23").Value = 0.08472035715708043
$ws.Range("G23").Value = 0.1795991340892121
